$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# This sheet is a daily/weekly price log for "Membrillo" (quince) sold
# at "Vega Modelo de Temuco". Two new daily records need to be inserted
# into the existing chronological block (rows 176-245), pushing the
# subsequent rows down. The new rows are NOT appended at the end; they
# are inserted in-place, which is why the diff shows nearly every row
# from 176 through 247 shifting by one (or two) positions.
# ---------------------------------------------------------------------

# 1) Insert the first new record at row 176 (pushes old rows 176-245
#    down to 177-246).
$ws.Rows(176).Insert()

$ws.Range("A176").Value = 10
$ws.Range("B176").Value = "Vega Modelo de Temuco"
$ws.Range("C176").Value = "La Araucanía"
$ws.Range("D176").Value = 45007
$ws.Range("E176").Value = 9
$ws.Range("F176").Value = "Fruta"
$ws.Range("G176").Value = 100104
$ws.Range("H176").Value = "Frutos de pepita"
$ws.Range("I176").Value = 100104003
$ws.Range("J176").Value = "Membrillo"
$ws.Range("K176").Value = "Champion"
$ws.Range("L176").Value = "Primera"
$ws.Range("M176").Value = 90
$ws.Range("N176").Value = 14000
$ws.Range("O176").Value = 14000
$ws.Range("P176").Value = 14000
$ws.Range("Q176").Value = "$/bandeja 18 kilos granel"
$ws.Range("R176").Value = "Región de O'Higgins"
$ws.Range("S176").Value = 778
$ws.Range("T176").Value = 18

# 2) Insert the second new record at row 241 (current numbering, i.e.
#    after the first insert already shifted things down by one). This
#    pushes the rows that are currently 241-246 (originally 240-245)
#    down to 242-247.
$ws.Rows(241).Insert()

$ws.Range("A241").Value = 10
$ws.Range("B241").Value = "Vega Modelo de Temuco"
$ws.Range("C241").Value = "La Araucanía"
$ws.Range("D241").Value = 45008
$ws.Range("E241").Value = 9
$ws.Range("F241").Value = "Fruta"
$ws.Range("G241").Value = 100104
$ws.Range("H241").Value = "Frutos de pepita"
$ws.Range("I241").Value = 100104003
$ws.Range("J241").Value = "Membrillo"
$ws.Range("K241").Value = "Champion"
$ws.Range("L241").Value = "Primera"
$ws.Range("M241").Value = 250
$ws.Range("N241").Value = 14000
$ws.Range("O241").Value = 14000
$ws.Range("P241").Value = 14000
$ws.Range("Q241").Value = "$/bandeja 18 kilos granel"
$ws.Range("R241").Value = "Región de O'Higgins"
$ws.Range("S241").Value = 778
$ws.Range("T241").Value = 18
